$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$handedBack = "Handed back: in sync with en-US"
$handbackDateTime = "2016-01-25 05:54:25"

# Overview sheet: de-de column (C) now shows the handed-back status
# for both a.md.md (row 2) and b.md.md (row 3).
$overview.Range("C2").Value = $handedBack
$overview.Range("C3").Value = $handedBack

# de-de sheet: record the handback for a.md.md (row 2) and its
# dependent b.md.md (row 3) - both share the same handed-off/handed-back
# target (a.md.md) and handback package file.
$dede.Range("B2").Value = $handedBack
$dede.Range("E2").Value = "a.md.md"
$dede.Range("F2").Value = $dede.Range("C2").Value()
$dede.Range("G2").Value = $handbackDateTime

$dede.Range("B3").Value = $handedBack
$dede.Range("E3").Value = "a.md.md"
$dede.Range("F3").Value = $dede.Range("C2").Value()
$dede.Range("G3").Value = $handbackDateTime

# Match the hyperlink/underline style used by the other "file name" style
# columns (A/C) for the newly populated E/F cells.
$dede.Range("E2").Style = "HyperLink"
$dede.Range("F2").Style = "HyperLink"
$dede.Range("E3").Style = "HyperLink"
$dede.Range("F3").Style = "HyperLink"

# Add the real hyperlinks behind the new "Latest Target File" / "Latest
# Handback File" cells, pointing at the same targets as the existing
# a.md.md / handback-xlf hyperlinks already on the sheet.
$aMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/3f60fd6c442af540561e24700a885d5f512ef46a/e2e/a.md.md"
$handbackXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8fcdb761a14fa60d7477198933472b9b2ec18410/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf"

$dede.Hyperlinks.Add($dede.Range("E2"), $aMdUrl, "", "", "a.md.md")
$dede.Hyperlinks.Add($dede.Range("F2"), $handbackXlfUrl, "", "", "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("E3"), $aMdUrl, "", "", "a.md.md")
$dede.Hyperlinks.Add($dede.Range("F3"), $handbackXlfUrl, "", "", "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf")

$wb.Save()
